$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the existing "IP" header (H1) onto the two new header
# cells so I1/J1 pick up the same style (bold, centered, bordered) as the
# rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..37: I0 and IF numeric values
$i0 = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,4,5,1,1,1,1,1,1,4,3)
$if = @(4,4,3,3,3,5,5,7,6,5,6,5,6,7,5,7,5,6,6,5,6,5,5,6,5,3,6,7,4,4,6,5,5,3,6,4)

for ($r = 2; $r -le 37; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0[$idx]
    $ws.Cells.Item($r, 10).Value = $if[$idx]
}
